# Update the wallet comments workbook:
#  - append three new challenge rows (#11, #12, #13) to "Designation Challenges"
#  - remove the "Q&A" worksheet entirely
#  - leave "Designation Requirements" as-is (it just shifts up one tab slot)

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false | Out-Null

$ws = $wb.Worksheets.Item("Designation Challenges")

# --- Row 12: "#11" / Receive - multi account alias issue ---
$ws.Cells.Item(12, 1).Value = "#11"
$ws.Cells.Item(12, 2).Value = "The account seems has alias, but I didn't find how to add alias to account(and this issue exist in other page that show multi account)"
$ws.Cells.Item(12, 3).Value = "Receive"
$ws.Cells.Item(12, 4).Value = "2a. Receive: Multi Account"
$ws.Cells.Item(12, 5).Value = "Receive"
$ws.Rows.Item(12).RowHeight = 27.6

# --- Row 13: "#12" / Send - missing "." key on keyboard ---
$ws.Cells.Item(13, 1).Value = "#12"
$ws.Cells.Item(13, 2).Value = "The keyboard lack the key "".""" 
$ws.Cells.Item(13, 3).Value = "Send"
$ws.Cells.Item(13, 4).Value = "3c. Send: Enter Amount"
$ws.Cells.Item(13, 5).Value = "Send"

# --- Row 14: "#13" / Backup wallet - mnemonic / HD wallet comment ---
$ws.Cells.Item(14, 1).Value = "#13"
$ws.Cells.Item(14, 2).Value = "I see the mnemonic words when backup and restore wallet, this means we should implement HD wallet(Hierarchical Deterministic Wallet), but it's hard and need lots of research. Even we get it done, the wallet key generated in mobile wallet can not be imported into coda client, so I suggest we manage the keys independently, which as the coda client do now."
$ws.Cells.Item(14, 3).Value = "Backup wallet"
$ws.Cells.Item(14, 4).Value = "0. Wallet Backup"
$ws.Cells.Item(14, 5).Value = "Backup wallet"
$ws.Rows.Item(14).RowHeight = 82.8

# Remove the Q&A sheet now that its open comments have been folded into the
# challenges list above.
$qa = $wb.Worksheets.Item("Q&A")
$qa.Delete() | Out-Null

# Make "Designation Challenges" the active tab, scrolled/selected near the
# newly added rows.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("B14").Select() | Out-Null
